# Apply updated crypto price/volume data (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.909.83"
$ws.Range("E2").Value = "  +4.10%  "
$ws.Range("D3").Value = "1.876.13"
$ws.Range("E3").Value = "  +3.53%  "
$ws.Range("E4").Value = "  +0.18%  "
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.97"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = "  +3.21%  "
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = "  +3.51%  "
$ws.Range("E7").Value = "  +0.09%  "
$__style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.49"
$ws.Range("D8").Style = $__style
$ws.Range("E8").Value = "  +10.89%  "
$ws.Range("E9").Value = "  +7.98%  "
$ws.Range("E10").Value = "  +3.78%  "
$ws.Range("E11").Value = "  +3.91%  "
$ws.Range("D12").Value = "2.144.64"
$ws.Range("E12").Value = "  +3.39%  "
$__style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.73"
$ws.Range("D13").Style = $__style
$ws.Range("E13").Value = "  +4.86%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$__style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.689"
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = "  +8.84%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.869.90"
$ws.Range("E15").Value = "  +2.40%  "
$__style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.81"
$ws.Range("D16").Style = $__style
$ws.Range("E16").Value = "  +8.75%  "
$ws.Range("D17").Value = "35.876.58"
$ws.Range("E17").Value = "  +4.09%  "
$__style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.84"
$ws.Range("D18").Style = $__style
$ws.Range("E18").Value = "  +3.57%  "
$ws.Range("D19").Value = "0.0₃0809"
$ws.Range("E19").Value = "  +4.49%  "
$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "248.24"
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = "  +2.36%  "
$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.51"
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = "  +11.33%  "
$__style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.83"
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = "  +17.13%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("E24").Value = "  +1.30%  "
$__style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.41"
$ws.Range("D25").Style = $__style
$ws.Range("E25").Value = "  +0.74%  "
$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.11"
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = "  +3.67%  "
$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.05"
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = "  +2.53%  "
$ws.Range("E28").Value = "  +2.45%  "
$ws.Range("E29").Value = "  +18.60%  "
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").Value = "3.309.84"
$ws.Range("E31").Value = "  +36.23%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.98"
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = "  +5.00%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$__style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0550"
$ws.Range("D33").Style = $__style
$ws.Range("E33").Value = "  +6.20%  "
$__style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.10"
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = "  +6.76%  "
$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.93"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = "  +5.26%  "
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "99.22"
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = "  +21.40%  "
$__style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.693"
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = "  +7.06%  "
$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.54"
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = "  +8.28%  "
$ws.Range("D39").Value = "1.367.51"
$ws.Range("E39").Value = "  +0.23%  "
$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.10"
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = "  +3.59%  "
$ws.Range("E41").Value = "  +6.07%  "
$ws.Range("E42").Value = "  +8.65%  "
$__style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.15"
$ws.Range("D43").Style = $__style
$ws.Range("E43").Value = "  +9.77%  "
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.27"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = "  +3.42%  "
$__style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.49"
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = "  +2.06%  "
$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.84"
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = "  +1.60%  "
$__style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.35"
$ws.Range("D47").Style = $__style
$ws.Range("E47").Value = "  +10.06%  "
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("D49").Value = "2.042.21"
$ws.Range("E49").Value = "  +3.40%  "
$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.67"
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = "  +3.55%  "
$ws.Range("E51").Value = "  +0.22%  "
